$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.148.80'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '2.470.41'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '582.63'
$ws.Range("D6").Value = '174.36'
$ws.Range("E6").Value = '  +3.27%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '2.916.60'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '25.38'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = '67.049.28'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '2.471.57'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("D20").Value = '348.53'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D23").Value = '69.46'
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.20'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = '2.596.17'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").Value = '0.0₃0898'
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("D30").Value = '499.29'
$ws.Range("D31").Value = '7.74'
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("D33").Value = '1.75'
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("D36").Value = '161.82'
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("D38").Value = '18.17'
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").Value = '142.68'
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = '0.581'
$ws.Range("E51").Value = '  -0.08%  '
